# Apply updated cryptos list values (price & volume change %) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.832.52'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '3.489.28'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '592.85'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '171.63'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.48%  '
$ws.Range("E7").Value = '  +0.00%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.589'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("E9").Value = '  +3.51%  '
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("D12").Value = '4.093.39'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("E13").Value = '  -0.49%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '28.93'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.56%  '
$ws.Range("D15").Value = '66.829.60'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = '3.494.66'
$ws.Range("E17").Value = '  +0.01%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '6.27'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.14%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '14.01'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.60%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '392.81'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.24%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '7.96'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.22%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '72.95'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("E25").Value = '  -1.21%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '10.17'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("E27").Value = '  -0.87%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.20%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '6.17'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.86%  '
$ws.Range("E30").Value = '  -3.51%  '
$ws.Range("E31").Value = '  -0.68%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '23.66'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("E33").Value = '  -1.30%  '
$ws.Range("E34").Value = '  -0.87%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '162.51'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +0.27%  '
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("E37").Value = '  -1.68%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '6.95'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +3.18%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '4.64'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.83%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0740'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.87%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '27.11'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").Value = '2.815.17'
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("E43").Value = '  -1.29%  '
$ws.Range("E45").Value = '  +2.01%  '
$ws.Range("E46").Value = '  -3.55%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '336.90'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -3.48%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '34.42'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.06%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.07'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.54%  '
$ws.Range("E50").Value = '  -1.12%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.40'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.39%  '
